# Fix binding student card
# Append the three missing students to Sheet1 (A4:B6) and leave the
# final selection where the author's Excel session left it (E13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 20120509
$ws.Range("B4").Value = "Nguyễn Đăng Khoa"

$ws.Range("A5").Value = 20120507
$ws.Range("B5").Value = "Võ Nhất Khanh"

$ws.Range("A6").Value = 20120511
$ws.Range("B6").Value = "Nguyễn Quốc Khoa"

# Match the saved cursor position recorded in the workbook
$ws.Range("E13").Select() | Out-Null
